# Atualização de bases das ligas, do dia: 20-02-2024 às 23:00
#
# This update (1) fixes the ordering of three pairs of fixtures that were
# recorded with swapped data rows, and (2) refreshes the odds snapshot
# for a handful of upcoming (not-yet-played) fixtures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the full data (columns B:AC - everything except the running row
#    index in column A) between mismatched row pairs.
# ---------------------------------------------------------------------------
$rowPairs = @(
    @(125, 126),
    @(133, 134),
    @(135, 137)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1" + ":AC$r1")
    $range2 = $ws.Range("B$r2" + ":AC$r2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}

# ---------------------------------------------------------------------------
# 2) Refresh odds snapshot values for upcoming fixtures.
# ---------------------------------------------------------------------------
$ws.Range("N307").Value2 = 1.3
$ws.Range("O307").Value2 = 5.5
$ws.Range("P307").Value2 = 9.5
$ws.Range("R307").Value2 = 1.875
$ws.Range("S307").Value2 = 1.975
$ws.Range("U307").Value2 = 1.825
$ws.Range("V307").Value2 = 2.025

$ws.Range("N308").Value2 = 1.571
$ws.Range("R308").Value2 = 2.05
$ws.Range("S308").Value2 = 1.8
$ws.Range("T308").Value2 = 2.5
$ws.Range("U308").Value2 = 2.025
$ws.Range("V308").Value2 = 1.825

$ws.Range("N310").Value2 = 1.571
$ws.Range("P310").Value2 = 5.5
$ws.Range("R310").Value2 = 1.775
$ws.Range("S310").Value2 = 2.1

$ws.Range("N311").Value2 = 2.1
$ws.Range("O311").Value2 = 3.25
$ws.Range("P311").Value2 = 3.6
$ws.Range("R311").Value2 = 1.85
$ws.Range("S311").Value2 = 2

$ws.Range("N312").Value2 = 1.8
$ws.Range("O312").Value2 = 3.4
$ws.Range("P312").Value2 = 4.75
$ws.Range("Q312").Value2 = -0.5
$ws.Range("R312").Value2 = 1.825
$ws.Range("S312").Value2 = 2.025
$ws.Range("U312").Value2 = 1.8
$ws.Range("V312").Value2 = 2.05

$ws.Range("R313").Value2 = 1.95
$ws.Range("S313").Value2 = 1.9
